$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "64.061.00"
$ws.Range("E2").Value = "  +1.73%  "

$ws.Range("D3").Value = "2.773.71"
$ws.Range("E3").Value = "  +2.62%  "

$ws.Range("E4").Value = "  +0.00%  "

Set-TextValue $ws.Range("D5") "585.96"
$ws.Range("E5").Value = "  +0.18%  "

Set-TextValue $ws.Range("D6") "161.28"
$ws.Range("E6").Value = "  +8.29%  "

Set-TextValue $ws.Range("D7") "0.620"
$ws.Range("E7").Value = "  +2.03%  "

$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("D9").Value = "2.782.77"
$ws.Range("E9").Value = "  +1.79%  "

$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("E11").Value = "  +0.68%  "

Set-TextValue $ws.Range("D12") "0.399"
$ws.Range("E12").Value = "  +2.82%  "

$ws.Range("E13").Value = "  +1.05%  "

$ws.Range("D14").Value = "3.274.95"
$ws.Range("E14").Value = "  +2.93%  "

Set-TextValue $ws.Range("D15") "27.59"
$ws.Range("E15").Value = "  +4.11%  "

$ws.Range("D16").Value = "63.957.97"
$ws.Range("E16").Value = "  +1.72%  "

Set-TextValue $ws.Range("D17") "0.0000159"
$ws.Range("E17").Value = "  +5.74%  "

$ws.Range("D18").Value = "2.782.89"
$ws.Range("E18").Value = "  +2.51%  "

Set-TextValue $ws.Range("D19") "12.41"
$ws.Range("E19").Value = "  +3.97%  "

Set-TextValue $ws.Range("D20") "5.04"
$ws.Range("E20").Value = "  +3.33%  "

Set-TextValue $ws.Range("D21") "367.62"
$ws.Range("E21").Value = "  +1.18%  "

Set-TextValue $ws.Range("D22") "7.07"
$ws.Range("E22").Value = "  +1.00%  "

Set-TextValue $ws.Range("D23") "0.572"
$ws.Range("E23").Value = "  +7.63%  "

$ws.Range("E24").Value = "  +0.76%  "

Set-TextValue $ws.Range("D25") "67.42"
$ws.Range("E25").Value = "  +3.17%  "

$ws.Range("E26").Value = "  +6.12%  "

Set-TextValue $ws.Range("D27") "8.78"
$ws.Range("E27").Value = "  +1.98%  "

$ws.Range("D28").Value = "0.0₃0972"
$ws.Range("E28").Value = "  +14.23%  "

$ws.Range("E29").Value = "  +0.42%  "

$ws.Range("E30").Value = "  +0.63%  "

Set-TextValue $ws.Range("D31") "7.28"
$ws.Range("E31").Value = "  +2.25%  "

Set-TextValue $ws.Range("D32") "1.27"
$ws.Range("E32").Value = "  +8.34%  "

Set-TextValue $ws.Range("D33") "172.65"
$ws.Range("E33").Value = "  +1.51%  "

Set-TextValue $ws.Range("D34") "5.09"
$ws.Range("E34").Value = "  +7.09%  "

Set-TextValue $ws.Range("D35") "20.81"
$ws.Range("E35").Value = "  +1.30%  "

$ws.Range("E36").Value = "  +0.12%  "

$ws.Range("E37").Value = "  +5.40%  "

$ws.Range("E38").Value = "  +0.56%  "

$ws.Range("E39").Value = "  +1.18%  "

Set-TextValue $ws.Range("D40") "4.28"
$ws.Range("E40").Value = "  +0.59%  "

Set-TextValue $ws.Range("D41") "6.29"
$ws.Range("E41").Value = "  +11.02%  "

Set-TextValue $ws.Range("D42") "341.41"
$ws.Range("E42").Value = "  -2.53%  "

Set-TextValue $ws.Range("D43") "40.03"
$ws.Range("E43").Value = "  +2.21%  "

Set-TextValue $ws.Range("D44") "22.46"
$ws.Range("E44").Value = "  +4.12%  "

Set-TextValue $ws.Range("D45") "22.71"
$ws.Range("E45").Value = "  +5.11%  "

Set-TextValue $ws.Range("D46") "0.0613"
$ws.Range("E46").Value = "  +3.13%  "

Set-TextValue $ws.Range("D47") "0.651"
$ws.Range("E47").Value = "  +1.68%  "

Set-TextValue $ws.Range("D48") "0.0262"
$ws.Range("E48").Value = "  +0.64%  "

Set-TextValue $ws.Range("D49") "138.69"
$ws.Range("E49").Value = "  +1.06%  "

$ws.Range("E50").Value = "  +1.83%  "

$ws.Range("D51").Value = "2.175.71"
$ws.Range("E51").Value = "  +1.84%  "
